$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.840987230004494
$ws.Range("C2").Value = 5.729659689673432
$ws.Range("D2").Value = 11.36463759379619
$ws.Range("F2").Value = 30.49735656803559
$ws.Range("G2").Value = 29.12426151942907
$ws.Range("H2").Value = 14.50101402468481
$ws.Range("J2").Value = 11.33628273765753
$ws.Range("K2").Value = 9.043114136576346
$ws.Range("O2").Value = 22.09889004117581
$ws.Range("B3").Value = 9.497542760162728
$ws.Range("C3").Value = 5.572449205096732
$ws.Range("D3").Value = 11.24997164847366
$ws.Range("F3").Value = 30.56731808119301
$ws.Range("G3").Value = 29.25942796862955
$ws.Range("H3").Value = 14.55689495268934
$ws.Range("J3").Value = 11.31648402814513
$ws.Range("K3").Value = 8.79086285541695
$ws.Range("O3").Value = 22.19740179909536
$ws.Range("B4").Value = 9.281111312286463
$ws.Range("C4").Value = 5.47323659175004
$ws.Range("D4").Value = 11.18119194356215
$ws.Range("F4").Value = 30.61814805911763
$ws.Range("G4").Value = 29.35197911113168
$ws.Range("H4").Value = 14.59351951457315
$ws.Range("J4").Value = 11.30656942361234
$ws.Range("K4").Value = 8.63285690076532
$ws.Range("O4").Value = 22.2626542117766
$ws.Range("B5").Value = 9.19164737421068
$ws.Range("C5").Value = 5.432176589371004
$ws.Range("D5").Value = 11.15360112595712
$ws.Range("F5").Value = 30.6408363375228
$ws.Range("G5").Value = 29.39208608819856
$ws.Range("H5").Value = 14.60902631938137
$ws.Range("J5").Value = 11.30309585856376
$ws.Range("K5").Value = 8.567772188076326
$ws.Range("O5").Value = 22.29044178879144
$ws.Range("B6").Value = 9.17671952743674
$ws.Range("C6").Value = 5.425321934781343
$ws.Range("D6").Value = 11.149046935773
$ws.Range("F6").Value = 30.64472281110752
$ws.Range("G6").Value = 29.39888991623224
$ws.Range("H6").Value = 14.61163637183759
$ws.Range("J6").Value = 11.30255338220758
$ws.Range("K6").Value = 8.556925606253161
$ws.Range("O6").Value = 22.29512812223449
$ws.Range("B7").Value = 9.27990971623298
$ws.Range("C7").Value = 5.472685331319568
$ws.Range("D7").Value = 11.18081803556224
$ws.Range("F7").Value = 30.61844605314283
$ws.Range("G7").Value = 29.35251034018491
$ws.Range("H7").Value = 14.59372628792092
$ws.Range("J7").Value = 11.30652027968439
$ws.Range("K7").Value = 8.6319818414708
$ws.Range("O7").Value = 22.26302412166736
$ws.Range("B8").Value = 9.723795038085905
$ws.Range("C8").Value = 5.676036722922609
$ws.Range("D8").Value = 11.32478082528106
$ws.Range("F8").Value = 30.51984224819379
$ws.Range("G8").Value = 29.1688760043293
$ws.Range("H8").Value = 14.51980177039843
$ws.Range("J8").Value = 11.3289927572695
$ws.Range("K8").Value = 8.956837684432275
$ws.Range("O8").Value = 22.13186672707338
$ws.Range("B9").Value = 10.5450629351983
$ws.Range("C9").Value = 6.051697541207584
$ws.Range("D9").Value = 11.61869485934874
$ws.Range("F9").Value = 30.38915983501574
$ws.Range("G9").Value = 28.8851289286194
$ws.Range("H9").Value = 14.39318018858559
$ws.Range("J9").Value = 11.39070005472781
$ws.Range("K9").Value = 9.56548898962904
$ws.Range("O9").Value = 21.91255828814751
$ws.Range("B10").Value = 11.11245244262921
$ws.Range("C10").Value = 6.311421914654656
$ws.Range("D10").Value = 11.83993636351362
$ws.Range("F10").Value = 30.33160846461652
$ws.Range("G10").Value = 28.72389591632655
$ws.Range("H10").Value = 14.31131773269755
$ws.Range("J10").Value = 11.44656220997649
$ws.Range("K10").Value = 9.990958232495764
$ws.Range("O10").Value = 21.77463612920275
$ws.Range("B11").Value = 11.36171092750149
$ws.Range("C11").Value = 6.425647116763502
$ws.Range("D11").Value = 11.94135451118835
$ws.Range("F11").Value = 30.31381866473346
$ws.Range("G11").Value = 28.66094216241467
$ws.Range("H11").Value = 14.27649788495714
$ws.Range("J11").Value = 11.47420336801226
$ws.Range("K11").Value = 10.17899649106185
$ws.Range("O11").Value = 21.71695306833856
$ws.Range("B12").Value = 11.45475256768624
$ws.Range("C12").Value = 6.468307930459352
$ws.Range("D12").Value = 11.97983827843622
$ws.Range("F12").Value = 30.30829081808279
$ws.Range("G12").Value = 28.63860826234678
$ws.Range("H12").Value = 14.26366032691151
$ws.Range("J12").Value = 11.4849853281963
$ws.Range("K12").Value = 10.24935043860299
$ws.Range("O12").Value = 21.69583939278024
$ws.Range("B13").Value = 11.43477538882773
$ws.Range("C13").Value = 6.459147007133362
$ws.Range("D13").Value = 11.97154711393263
$ws.Range("F13").Value = 30.30942755130877
$ws.Range("G13").Value = 28.64335114874508
$ws.Range("H13").Value = 14.26640964880789
$ws.Range("J13").Value = 11.48264933212765
$ws.Range("K13").Value = 10.2342372422635
$ws.Range("O13").Value = 21.70035411331987
$ws.Range("B14").Value = 11.36939291195395
$ws.Range("C14").Value = 6.429168906087599
$ws.Range("D14").Value = 11.94451919698417
$ws.Range("F14").Value = 30.3133396499889
$ws.Range("G14").Value = 28.65907450919026
$ws.Range("H14").Value = 14.27543475757841
$ws.Range("J14").Value = 11.47508413035465
$ws.Range("K14").Value = 10.18480195073238
$ws.Range("O14").Value = 21.71520139737479
$ws.Range("B15").Value = 11.32916676607343
$ws.Range("C15").Value = 6.41072830168051
$ws.Range("D15").Value = 11.9279731371109
$ws.Range("F15").Value = 30.31589339319336
$ws.Range("G15").Value = 28.66890187522228
$ws.Range("H15").Value = 14.28100821259951
$ws.Range("J15").Value = 11.47049105279057
$ws.Range("K15").Value = 10.15440877031044
$ws.Range("O15").Value = 21.72439087627915
$ws.Range("B16").Value = 11.09597814654317
$ws.Range("C16").Value = 6.303875470735937
$ws.Range("D16").Value = 11.83332125187613
$ws.Range("F16").Value = 30.33294011548232
$ws.Range("G16").Value = 28.72822015446085
$ws.Range("H16").Value = 14.31364198356876
$ws.Range("J16").Value = 11.44480017442227
$ws.Range("K16").Value = 9.978553376057461
$ws.Range("O16").Value = 21.77850783358522
$ws.Range("B17").Value = 10.95060177976748
$ws.Range("C17").Value = 6.237297622201798
$ws.Range("D17").Value = 11.77543055166407
$ws.Range("F17").Value = 30.34554856271795
$ws.Range("G17").Value = 28.76727959446471
$ws.Range("H17").Value = 14.33428150568202
$ws.Range("J17").Value = 11.42960667753961
$ws.Range("K17").Value = 9.869216197534874
$ws.Range("O17").Value = 21.81300401642876
$ws.Range("B18").Value = 10.86615597963914
$ws.Range("C18").Value = 6.19863595587884
$ws.Range("D18").Value = 11.74220816036059
$ws.Range("F18").Value = 30.35359028533679
$ws.Range("G18").Value = 28.79072306641901
$ws.Range("H18").Value = 14.34638057418202
$ws.Range("J18").Value = 11.42107796948186
$ws.Range("K18").Value = 9.805812922159438
$ws.Range("O18").Value = 21.83332129539836
$ws.Range("B19").Value = 10.83742409961794
$ws.Range("C19").Value = 6.185483546935934
$ws.Range("D19").Value = 11.7309734699218
$ws.Range("F19").Value = 30.35644862320099
$ws.Range("G19").Value = 28.79882816969398
$ws.Range("H19").Value = 14.35051622945565
$ws.Range("J19").Value = 11.41822656254859
$ws.Range("K19").Value = 9.784259055525881
$ws.Range("O19").Value = 21.84028205974187
$ws.Range("B20").Value = 10.96616372663613
$ws.Range("C20").Value = 6.244423241335943
$ws.Range("D20").Value = 11.78158560886643
$ws.Range("F20").Value = 30.34412462504903
$ws.Range("G20").Value = 28.76302040010761
$ws.Range("H20").Value = 14.33206082093967
$ws.Range("J20").Value = 11.43120233736603
$ws.Range("K20").Value = 9.880909122539231
$ws.Range("O20").Value = 21.80928255970866
$ws.Range("B21").Value = 11.38863443941172
$ws.Range("C21").Value = 6.437990536226773
$ws.Range("D21").Value = 11.95245606657579
$ws.Range("F21").Value = 30.31215775230349
$ws.Range("G21").Value = 28.65441523831035
$ws.Range("H21").Value = 14.27277442066693
$ws.Range("J21").Value = 11.47729771594691
$ws.Range("K21").Value = 10.19934586611161
$ws.Range("O21").Value = 21.71082056970791
$ws.Range("B22").Value = 11.65686310983796
$ws.Range("C22").Value = 6.561025281175621
$ws.Range("D22").Value = 12.06457496888811
$ws.Range("F22").Value = 30.29831166060227
$ws.Range("G22").Value = 28.59221435489446
$ws.Range("H22").Value = 14.23605553642933
$ws.Range("J22").Value = 11.50925598517819
$ws.Range("K22").Value = 10.40247256982101
$ws.Range("O22").Value = 21.6507239649054
$ws.Range("B23").Value = 11.5144477438416
$ws.Range("C23").Value = 6.495686037678572
$ws.Range("D23").Value = 12.0047048959332
$ws.Range("F23").Value = 30.30505632219224
$ws.Range("G23").Value = 28.62460538355543
$ws.Range("H23").Value = 14.25546751036216
$ws.Range("J23").Value = 11.49203359151003
$ws.Range("K23").Value = 10.29453486670986
$ws.Range("O23").Value = 21.68240866588832
$ws.Range("B24").Value = 10.95913086518272
$ws.Range("C24").Value = 6.241202946411961
$ws.Range("D24").Value = 11.77880271894053
$ws.Range("F24").Value = 30.34476591742337
$ws.Range("G24").Value = 28.76494290614207
$ws.Range("H24").Value = 14.33306406630826
$ws.Range("J24").Value = 11.43048029673335
$ws.Range("K24").Value = 9.875624440594544
$ws.Range("O24").Value = 21.81096351981383
$ws.Range("B25").Value = 10.32881930346308
$ws.Range("C25").Value = 5.952783931796669
$ws.Range("D25").Value = 11.53812258712521
$ws.Range("F25").Value = 30.41777429167761
$ws.Range("G25").Value = 28.9536429365049
$ws.Range("H25").Value = 14.42547260869095
$ws.Range("J25").Value = 11.37213991790844
$ws.Range("K25").Value = 9.56548898962904
$ws.Range("O25").Value = 21.96781983487079
